# Add a new slide (12) using the "Title and Content" layout (same layout
# used by every other content slide in this deck) and populate its title
# and body placeholders, matching the new "Biological implications of
# small-world properties" slide added to the deck.

$p = $ppt.ActivePresentation

$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Biological implications of small-world properties"

# Body / content placeholder - one paragraph per bullet line
$bodyLines = @(
    "Neural Systems: Balancing Specialization and Integration",
    "Genetic and Metabolic Networks: Efficiency and Evolution",
    "Plant Communication: Underground Signaling Networks",
    "Cellular Communication: Coordinated Responses",
    "Ecological Systems: Synchronization and Stability"
)

$bodyRange = $s.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = $bodyLines[0]
for ($i = 1; $i -lt $bodyLines.Count; $i++) {
    [void]$bodyRange.InsertAfter("`r" + $bodyLines[$i])
}
